# Apply weekly update to "Fruta, Comercializadora del Agro de Limarí - Palta"
# The new week's data (date 2022-09-05 / serial 44826) is inserted as 5 new rows
# right before the existing row 479, pushing all subsequent rows down by 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows at row 479 (existing rows 479:517 shift down to 484:522)
$ws.Rows("479:483").Insert()

# Common (constant) values shared by every row in this data block
$commonA = 2
$commonB = "Comercializadora del Agro de Limarí"
$commonC = "Coquimbo"
$commonE = 4
$commonF = "Fruta"
$commonG = 100106
$commonH = "Oleaginosos"
$commonI = 100106002
$commonJ = "Palta"
$commonR = "Provincia de Limarí"
$commonT = 1

# New-row data: row -> (K variedad, L calidad, M volumen, N precio min, O precio max, P precio prom, Q unidad, S precio/kg)
$newRows = @(
    @{ Row = 479; D = 44826; K = "Edranol"; L = "Primera";      M = 500; N = 1400; O = 1500; P = 1450; Q = "`$/kilo (en caja de 17 kilos)"; S = 1450 },
    @{ Row = 480; D = 44826; K = "Edranol"; L = "Segunda";      M = 360; N = 1000; O = 1100; P = 1050; Q = "`$/kilo (en caja de 17 kilos)"; S = 1050 },
    @{ Row = 481; D = 44826; K = "Hass";    L = "1a nueva(o)";  M = 540; N = 1800; O = 1900; P = 1850; Q = "`$/kilo (en caja de 17 kilos)"; S = 1850 },
    @{ Row = 482; D = 44826; K = "Hass";    L = "2a nueva(o)";  M = 400; N = 1600; O = 1700; P = 1650; Q = "`$/kilo (en caja de 17 kilos)"; S = 1650 },
    @{ Row = 483; D = 44826; K = "Hass";    L = "3a nueva (o)"; M = 300; N = 1300; O = 1400; P = 1350; Q = "`$/kilo (en caja de 17 kilos)"; S = 1350 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2 = $commonA
    $ws.Cells.Item($row, 2).Value = $commonB
    $ws.Cells.Item($row, 3).Value = $commonC
    $ws.Cells.Item($row, 4).Value2 = $r.D
    $ws.Cells.Item($row, 5).Value2 = $commonE
    $ws.Cells.Item($row, 6).Value = $commonF
    $ws.Cells.Item($row, 7).Value2 = $commonG
    $ws.Cells.Item($row, 8).Value = $commonH
    $ws.Cells.Item($row, 9).Value2 = $commonI
    $ws.Cells.Item($row, 10).Value = $commonJ
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $r.N
    $ws.Cells.Item($row, 15).Value2 = $r.O
    $ws.Cells.Item($row, 16).Value2 = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $commonR
    $ws.Cells.Item($row, 19).Value2 = $r.S
    $ws.Cells.Item($row, 20).Value2 = $commonT
}


